# Adds two new columns, I ("I0") and J ("IF"), to the sheet, mirroring the
# header style already used by the other header cells (bold + border, style
# index 1 in the original workbook) and filling in the per-row numeric
# values for rows 2-70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they pick up the same bold/border/center style used by
# the rest of row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Per-row data for columns I and J, rows 2-70 (69 values each).
$iValues = @(6,8,5,5,6,6,7,6,5,10,8,6,7,9,6,8,6,6,8,7,7,8,6,6,1,7,7,7,6,8,5,7,8,10,5,10,9,8,6,5,8,8,7,7,5,7,5,8,7,5,6,6,10,5,8,7,8,9,4,7,8,8,7,6,7,7,8,3,2)
$jValues = @(7,8,5,6,7,7,8,7,5,10,9,6,7,9,7,8,7,8,8,8,7,8,6,7,2,7,7,7,6,9,5,7,8,10,6,10,9,9,7,6,8,8,7,7,6,8,5,8,8,5,6,7,10,6,8,8,8,9,5,7,8,8,7,7,7,7,8,3,2)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

Write-Output "Added columns I (I0) and J (IF) for rows 1-70"
